$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Shift the whole table one column to the right and one row down by
#    inserting a new blank column A and a new blank row 1. Excel moves
#    values, number formats, fills, and column widths automatically.
# ---------------------------------------------------------------------
$ws.Columns("A:A").Insert()
$ws.Rows("1:1").Insert()

# ---------------------------------------------------------------------
# 2. Refresh the clustering results: update the cells whose values
#    changed (re-run of the clustering produced a new table).
# ---------------------------------------------------------------------
$ws.Range("F3").Value = 2.7
$ws.Range("J3").Value = 4.0999999999999996
$ws.Range("K3").Value = 44.7

$ws.Range("C4").Value = 6
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 4.8
$ws.Range("G4").Value = 5.8
$ws.Range("H4").Value = 19.2
$ws.Range("I4").Value = 2.5
$ws.Range("J4").Value = 7.8
$ws.Range("K4").Value = 18.5
$ws.Range("L4").Value = 16
$ws.Range("M4").Value = 25.2

$ws.Range("C5").Value = 6
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 4.2
$ws.Range("G5").Value = 1.3
$ws.Range("H5").Value = 22
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 24.2
$ws.Range("L5").Value = 33.200000000000003
$ws.Range("M5").Value = 7.7

$ws.Range("J9").Value = 7

# ---------------------------------------------------------------------
# 3. The shading (highlight) that tracks the row maximum moved along
#    with the refreshed numbers: M4 is now the row-4 maximum (gets the
#    gray highlight) and M5 no longer is (highlight removed). Reuse the
#    existing highlighted / plain styles already on the sheet instead of
#    inventing new ones, by copying formats with Paste Special.
# ---------------------------------------------------------------------
$ws.Range("K3").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C3").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Refresh the view: hide gridlines, reset zoom to 100%, scroll so
#    column C is left-most, and leave the selection on F5 like in the
#    saved workbook.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.DisplayGridlines = $false
$win.Zoom = 100
$ws.Range("F5").Select()
$win.ScrollColumn = 3
$win.ScrollRow = 1
